$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text so numeric-looking strings
# (e.g. "1.10", "0.120") round-trip exactly instead of being coerced
# to numbers by Excel's smart type inference.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "36.235.98"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "2.003.39"
$ws.Range("E3").Value = "  +5.54%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "243.88"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").Value = "0.658"
$ws.Range("E6").Value = "  -5.24%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "44.25"
$ws.Range("E8").Value = "  +2.31%  "
$ws.Range("D9").Value = "61.51"
$ws.Range("E9").Value = "  +6.87%  "
$ws.Range("D10").Value = "0.362"
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("D11").Value = "0.0710"
$ws.Range("E11").Value = "  -6.16%  "
$ws.Range("D12").Value = "0.0977"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "14.22"
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("D14").Value = "2.295.77"
$ws.Range("E14").Value = "  +5.64%  "
$ws.Range("D15").Value = "0.801"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "2.006.77"
$ws.Range("E16").Value = "  +5.17%  "
$ws.Range("D17").Value = "4.85"
$ws.Range("E17").Value = "  -3.97%  "
$ws.Range("D18").Value = "36.210.45"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").Value = "70.82"
$ws.Range("E19").Value = "  -4.18%  "
$ws.Range("D20").Value = "0.0₃0807"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").Value = "235.52"
$ws.Range("E21").Value = "  -4.76%  "
$ws.Range("D22").Value = "12.68"
$ws.Range("E22").Value = "  -2.76%  "
$ws.Range("D23").Value = "4.85"
$ws.Range("E23").Value = "  -7.05%  "
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  -9.51%  "
$ws.Range("D26").Value = "165.52"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").Value = "8.56"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "19.26"
$ws.Range("E28").Value = "  +4.58%  "
$ws.Range("E29").Value = "  -10.70%  "
$ws.Range("D30").Value = "0.120"
$ws.Range("E30").Value = "  -6.27%  "
$ws.Range("D31").Value = "21.15"
$ws.Range("E31").Value = "  +48.01%  "
$ws.Range("D32").Value = "4.31"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "0.0578"
$ws.Range("E33").Value = "  -4.03%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "1.87"
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.0860"
$ws.Range("E36").Value = "  +17.56%  "
$ws.Range("D37").Value = "3.94"
$ws.Range("E37").Value = "  -7.64%  "
$ws.Range("D38").Value = "2.08"
$ws.Range("E38").Value = "  +5.48%  "
$ws.Range("D39").Value = "0.841"
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D40").Value = "1.31"
$ws.Range("E40").Value = "  -11.62%  "
$ws.Range("D41").Value = "0.0213"
$ws.Range("E41").Value = "  -6.14%  "
$ws.Range("D42").Value = "94.85"
$ws.Range("E42").Value = "  -4.92%  "
$ws.Range("D43").Value = "1.10"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").Value = "2.75"
$ws.Range("E44").Value = "  +15.25%  "
$ws.Range("D45").Value = "15.86"
$ws.Range("E45").Value = "  -7.64%  "
$ws.Range("D46").Value = "1.300.60"
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("D47").Value = "0.0809"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "2.76"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").Value = "2.18"
$ws.Range("E49").Value = "  -8.38%  "
$ws.Range("D50").Value = "2.171.80"
$ws.Range("E50").Value = "  +4.73%  "
$ws.Range("E51").Value = "  +14.54%  "

# Restore the original (default/no explicit) style so the cells keep
# matching their pre-edit appearance (no stray Text number-format left behind).
$dataRange.Style = "Normal"
